# Update "想去人数" (column F) values in the "展览" and "全部类型" sheets
# to reflect refreshed counts at the time the gh-pages output was regenerated.

$wb = $excel.ActiveWorkbook

# Row -> (old, new) value map for column F
$updates = @{
    3  = 10893
    4  = 331
    5  = 992
    6  = 216
    7  = 1355
    8  = 8361
    10 = 472
    11 = 617
    12 = 229
    13 = 140
    14 = 3349
    16 = 332
    18 = 848
    20 = 1082
    22 = 139
    23 = 1883
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
